$d = $word.ActiveDocument

# 1) "Supuestos y dependencias - Por decidir" -> "Supuestos y dependencias - Neifi"
$d.Content.Find.Execute("Supuestos y dependencias " + [char]0x2013 + " Por decidir", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Supuestos y dependencias " + [char]0x2013 + " Neifi", 2) | Out-Null

# 2) Move the _GoBack bookmark from the end of the "En caso de..." paragraph
#    to the end of the "Supuestos y dependencias - Neifi" heading paragraph
#    (right after the new text, inside the same paragraph).
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$heading = $d.Content.Find.Execute("Supuestos y dependencias " + [char]0x2013 + " Neifi") | Out-Null
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Supuestos y dependencias*Neifi*") {
        $target = $para
    }
}
$insertPos = $target.Range.End - 1
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3) "Propuesta economica - Por decidir " -> "Propuesta economica - Carlos"
$d.Content.Find.Execute("Propuesta econ" + [char]0x00F3 + "mica " + [char]0x2013 + " Por decidir ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Propuesta econ" + [char]0x00F3 + "mica " + [char]0x2013 + " Carlos", 2) | Out-Null
